$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.756.11'
$ws.Range('E2').Value = '  +4.71%  '
$ws.Range('D3').Value = '1.879.65'
$ws.Range('E3').Value = '  +2.88%  '
$ws.Range('E4').Value = '  -0.52%  '
$ws.Range('D5').Value = "'338.47"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.93%  '
$ws.Range('D6').Value = "'0.9996"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.50%  '
$ws.Range('D7').Value = "'0.4728"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.14%  '
$ws.Range('D8').Value = "'0.4035"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.65%  '
$ws.Range('D9').Value = "'47.66"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.77%  '
$ws.Range('D10').Value = "'0.08063"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E11').Value = '  +3.89%  '
$ws.Range('D12').Value = "'22.28"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.83%  '
$ws.Range('D13').Value = "'6.060"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.99%  '
$ws.Range('D14').Value = '1.874.20'
$ws.Range('E14').Value = '  +2.33%  '
$ws.Range('D15').Value = "'7.304"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.65%  '
$ws.Range('D16').Value = "'90.76"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('D17').Value = "'1.001"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('E18').Value = '  +1.91%  '
$ws.Range('D19').Value = "'0.06604"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.46%  '
$ws.Range('D20').Value = "'17.68"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.29%  '
$ws.Range('D22').Value = '28.780.50'
$ws.Range('E22').Value = '  +4.87%  '
$ws.Range('D23').Value = "'5.510"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.22%  '
$ws.Range('E24').Value = '  +2.34%  '
$ws.Range('D25').Value = "'2.261"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.79%  '
$ws.Range('D26').Value = '2.098.99'
$ws.Range('E26').Value = '  +2.51%  '
$ws.Range('D27').Value = "'160.47"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.39%  '
$ws.Range('D28').Value = "'19.87"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.71%  '
$ws.Range('D29').Value = "'2.133"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.51%  '
$ws.Range('D30').Value = "'5.508"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.38%  '
$ws.Range('D31').Value = "'120.10"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.62%  '
$ws.Range('D32').Value = "'0.9861"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.09%  '
$ws.Range('D33').Value = "'0.09560"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.53%  '
$ws.Range('D34').Value = "'3.657"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.35%  '
$ws.Range('D35').Value = "'1.396"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.28%  '
$ws.Range('D36').Value = "'5.380"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.61%  '
$ws.Range('D37').Value = "'0.06184"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.39%  '
$ws.Range('D38').Value = "'0.02279"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.41%  '
$ws.Range('D39').Value = "'8.524"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.17%  '
$ws.Range('D40').Value = "'1.186"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.32%  '
$ws.Range('D41').Value = "'0.5970"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.58%  '
$ws.Range('D42').Value = "'0.9991"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.38%  '
$ws.Range('D43').Value = "'0.1890"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.12%  '
$ws.Range('D44').Value = "'10.39"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.42%  '
$ws.Range('D45').Value = "'1.264"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = "'0.5598"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.74%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'12.18"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.22%  '
$ws.Range('E48').Value = '  +5.15%  '
$ws.Range('D49').Value = "'0.07235"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +9.42%  '
$ws.Range('D50').Value = "'2.108"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +13.70%  '
$ws.Range('D51').Value = "'112.47"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.79%  '
